# Add a new "EC_Prism_Template" worksheet derived from "EC_Template",
# with a trimmed-down column set (A:G) and two new headers pulled in
# from a Prism-flavoured variant of the report.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the template sheet (keeps fonts/borders/merges/col widths)
# and place the copy right after the original.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "EC_Prism_Template"

# Drop the "SMV No Tag File Name" column (F) -- this shifts what was
# G ("EC Exists") left into F.
$ws2.Range("F1").EntireColumn.Delete()

# Drop the old "Tags"/"No Tags" sub-header block (was H:S, now G:R
# after the first delete), leaving just A:F.
$ws2.Range("G1:R3").EntireColumn.Delete()

# Re-label the remaining "SMV Tag File Name" header for the Prism variant.
$ws2.Range("E2").Value = "Prism No Tag File Name"

# New trailing column with its own header.
$ws2.Range("G2").Value = "Prob for EC"

# Match formatting of the new column G to its neighbours: G1 follows the
# title-row style, G2/G3 follow the bordered header style used by the
# rest of row 2/3.
$ws2.Range("B3").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$ws2.Range("A3").Copy()
$ws2.Range("G2").PasteSpecial(-4122)
$ws2.Range("G3").PasteSpecial(-4122)
$ws2.Range("G2:G3").Merge()

# Row 3 no longer has the tall wrapped "Tags"/"No Tags" text, so its
# height shrinks; row 4 exists but stays blank.
$ws2.Rows.Item(3).RowHeight = 13.8
$ws2.Range("G4").Borders.LineStyle = 0
$ws2.Rows.Item(4).RowHeight = 13.8

$ws2.Range("A1").Select()
